$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.734.04'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +8.41%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.947.62'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +7.02%  '
$ws.Range("E4").Value = '  -0.26%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '341.55'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.99%  '
$ws.Range("E6").Value = '  -0.19%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4782'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.48%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4145'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +8.99%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '48.45'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.96%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08272'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.24%  '
$ws.Range("E11").Value = '  +8.66%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.66'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +7.86%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.941.49'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.36%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.194'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.24%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.431'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.21%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '92.31'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.38%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.002'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.14%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001064'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.34%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06675'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.31%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.08'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.62%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.000'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.28%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '29.688.72'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +8.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.606'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +6.23%  '
$ws.Range("E24").Value = '  +4.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.283'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.04%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.171.36'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.98%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '160.38'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.62%  '
$ws.Range("E28").Value = '  +4.73%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.196'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +7.36%  '
$ws.Range("E30").Value = '  +7.24%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '122.39'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.09%  '
$ws.Range("E32").Value = '  +10.12%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09640'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.74%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.473'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +11.88%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.680'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.26%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.487'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.25%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06314'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.89%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02329'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +6.94%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.606'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +6.23%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.200'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.23%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6109'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.32%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '10.69'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +7.73%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1902'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.84%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.000'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.23%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.286'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.05%  '
$ws.Range("B46").Value = 'RenderToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.357'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +31.14%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.54'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +6.42%  '
$ws.Range("B48").Value = 'Decentraland'
$ws.Range("C48").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5710'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.88%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.004'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +7.24%  '
$ws.Range("E50").Value = '  +12.72%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '114.01'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.94%  '

Write-Output "Applied all changes"
